$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Values are prefixed with a leading apostrophe so Excel stores them
# as literal text (matching the original inlineStr cells) instead of
# auto-converting number-like strings (e.g. "1.001") into numerics.

$ws.Range("D2").Value = "'29.891.11"
$ws.Range("E2").Value = "'  +0.17%  "
$ws.Range("D3").Value = "'1.891.66"
$ws.Range("E3").Value = "'  -0.05%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'0.7723"
$ws.Range("E5").Value = "'  -1.40%  "
$ws.Range("D6").Value = "'243.58"
$ws.Range("E6").Value = "'  -0.14%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("D8").Value = "'0.3124"
$ws.Range("E8").Value = "'  -0.62%  "
$ws.Range("D9").Value = "'25.64"
$ws.Range("E9").Value = "'  +1.35%  "
$ws.Range("D10").Value = "'0.07234"
$ws.Range("E10").Value = "'  +0.53%  "
$ws.Range("D11").Value = "'0.08717"
$ws.Range("E11").Value = "'  +7.72%  "
$ws.Range("D12").Value = "'2.117.34"
$ws.Range("E12").Value = "'  +12.75%  "
$ws.Range("D13").Value = "'0.7707"
$ws.Range("E13").Value = "'  +0.84%  "
$ws.Range("D14").Value = "'5.399"
$ws.Range("E14").Value = "'  -1.66%  "
$ws.Range("D15").Value = "'94.28"
$ws.Range("E15").Value = "'  +2.12%  "
$ws.Range("D16").Value = "'6.209"
$ws.Range("E16").Value = "'  +1.00%  "
$ws.Range("D17").Value = "'29.896.29"
$ws.Range("E17").Value = "'  +0.16%  "
$ws.Range("D18").Value = "'13.91"
$ws.Range("E18").Value = "'  -0.44%  "
$ws.Range("B19").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "'2.352.83"
$ws.Range("E19").Value = "'  +9.52%  "
$ws.Range("B20").Value = "'BitcoinCash"
$ws.Range("C20").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'245.17"
$ws.Range("E20").Value = "'  +0.72%  "
$ws.Range("D21").Value = "'0.000007874"
$ws.Range("E21").Value = "'  +1.05%  "
$ws.Range("D22").Value = "'8.173"
$ws.Range("E22").Value = "'  +0.33%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "'  -0.12%  "
$ws.Range("E24").Value = "'  -0.03%  "
$ws.Range("D25").Value = "'0.1590"
$ws.Range("E25").Value = "'  -3.11%  "
$ws.Range("D26").Value = "'9.512"
$ws.Range("E26").Value = "'  +0.93%  "
$ws.Range("D27").Value = "'162.46"
$ws.Range("E27").Value = "'  -0.36%  "
$ws.Range("D28").Value = "'18.81"
$ws.Range("E28").Value = "'  +0.43%  "
$ws.Range("D29").Value = "'2.045"
$ws.Range("E29").Value = "'  -0.38%  "
$ws.Range("E30").Value = "'  +1.41%  "
$ws.Range("E31").Value = "'  -0.28%  "
$ws.Range("D32").Value = "'4.515"
$ws.Range("E32").Value = "'  +0.39%  "
$ws.Range("D33").Value = "'4.120"
$ws.Range("E33").Value = "'  +0.09%  "
$ws.Range("D34").Value = "'0.05438"
$ws.Range("E34").Value = "'  -2.27%  "
$ws.Range("D35").Value = "'1.247"
$ws.Range("E35").Value = "'  -1.62%  "
$ws.Range("D36").Value = "'0.7511"
$ws.Range("E36").Value = "'  +0.95%  "
$ws.Range("D37").Value = "'1.006"
$ws.Range("E37").Value = "'  +0.75%  "
$ws.Range("D38").Value = "'2.695"
$ws.Range("E38").Value = "'  +2.99%  "
$ws.Range("D39").Value = "'0.01983"
$ws.Range("E39").Value = "'  +3.20%  "
$ws.Range("D40").Value = "'2.785"
$ws.Range("E40").Value = "'  +0.04%  "
$ws.Range("D41").Value = "'0.4514"
$ws.Range("E41").Value = "'  +2.15%  "
$ws.Range("B42").Value = "'Aave"
$ws.Range("C42").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'73.39"
$ws.Range("E42").Value = "'  -0.51%  "
$ws.Range("B43").Value = "'Maker"
$ws.Range("C43").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'1.099.13"
$ws.Range("E43").Value = "'  -4.32%  "
$ws.Range("D44").Value = "'6.092"
$ws.Range("E44").Value = "'  +3.94%  "
$ws.Range("D45").Value = "'0.8552"
$ws.Range("E45").Value = "'  +0.49%  "
$ws.Range("E46").Value = "'  -0.07%  "
$ws.Range("D47").Value = "'2.185.45"
$ws.Range("E47").Value = "'  +6.77%  "
$ws.Range("D48").Value = "'103.52"
$ws.Range("E48").Value = "'  -0.38%  "
$ws.Range("E49").Value = "'  +0.33%  "
$ws.Range("D50").Value = "'7.612"
$ws.Range("E50").Value = "'  +2.01%  "
$ws.Range("D51").Value = "'9.848"
$ws.Range("E51").Value = "'  -1.31%  "
